$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feriekasse")

$ws.Range("G2").Value = "Man Utd"
$ws.Range("C7").Value = "Brndby"
$ws.Range("E7").Value = "FCK"
$ws.Range("G7").Value = "FC Midtjylland"
